$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'57.199.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "'3.063.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.97%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'513.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.42%  "
$ws.Range("D6").Value = "'142.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.48%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +3.04%  "
$ws.Range("D9").Value = "'7.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("E10").Value = "  +3.72%  "
$ws.Range("E11").Value = "  +6.01%  "
$ws.Range("D12").Value = "'3.593.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.26%  "
$ws.Range("E13").Value = "  +3.11%  "
$ws.Range("D14").Value = "'25.74"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("E15").Value = "  +3.47%  "
$ws.Range("D16").Value = "'57.345.04"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "'3.078.58"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.58%  "
$ws.Range("D18").Value = "'6.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("D19").Value = "'12.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.81%  "
$ws.Range("D20").Value = "'8.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.06%  "
$ws.Range("D21").Value = "'334.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.66%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D24").Value = "'65.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.11%  "
$ws.Range("E25").Value = "  +7.01%  "
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'0.0₃0928"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.13%  "
$ws.Range("D28").Value = "'6.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.61%  "
$ws.Range("D29").Value = "'7.07"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.79%  "
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("D31").Value = "'20.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.11%  "
$ws.Range("E32").Value = "  +2.39%  "
$ws.Range("D33").Value = "'154.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.41%  "
$ws.Range("E34").Value = "  +3.50%  "
$ws.Range("E35").Value = "  +5.49%  "
$ws.Range("D36").Value = "'26.02"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.75%  "
$ws.Range("E37").Value = "  +4.22%  "
$ws.Range("D38").Value = "'0.0682"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("D39").Value = "'3.103.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.28%  "
$ws.Range("D40").Value = "'36.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.91%  "
$ws.Range("E41").Value = "  +3.86%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +5.08%  "
$ws.Range("D44").Value = "'2.264.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.19%  "
$ws.Range("E45").Value = "  +8.20%  "
$ws.Range("E46").Value = "  +4.26%  "
$ws.Range("E47").Value = "  +3.67%  "
$ws.Range("E48").Value = "  +9.22%  "
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("E50").Value = "  +4.56%  "
$ws.Range("E51").Value = "  +7.33%  "
